$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.620.07"
$ws.Range("E2").Value = "  +14.57%  "
$ws.Range("D3").Value = "1.798.84"
$ws.Range("E3").Value = "  +7.43%  "
$ws.Range("D4").Value = "0.993"
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("D5").Value = "231.97"
$ws.Range("E5").Value = "  +5.69%  "
$ws.Range("D6").Value = "0.550"
$ws.Range("E6").Value = "  +5.28%  "
$ws.Range("D7").Value = "0.991"
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("D8").Value = "31.48"
$ws.Range("E8").Value = "  +5.99%  "
$ws.Range("D9").Value = "46.35"
$ws.Range("E9").Value = "  +5.61%  "
$ws.Range("D10").Value = "0.284"
$ws.Range("E10").Value = "  +7.11%  "
$ws.Range("D11").Value = "0.0675"
$ws.Range("E11").Value = "  +9.10%  "
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("D13").Value = "2.052.63"
$ws.Range("E13").Value = "  +7.24%  "
$ws.Range("D14").Value = "1.803.65"
$ws.Range("E14").Value = "  +7.71%  "
$ws.Range("D15").Value = "0.640"
$ws.Range("E15").Value = "  +3.46%  "
$ws.Range("D16").Value = "34.536.54"
$ws.Range("E16").Value = "  +14.21%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "10.25"
$ws.Range("E17").Value = "  -4.49%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "4.33"
$ws.Range("E18").Value = "  +8.16%  "
$ws.Range("D19").Value = "70.25"
$ws.Range("E19").Value = "  +6.93%  "
$ws.Range("D20").Value = "263.94"
$ws.Range("E20").Value = "  +6.84%  "
$ws.Range("D21").Value = "0.0₃0758"
$ws.Range("E21").Value = "  +5.55%  "
$ws.Range("D22").Value = "0.991"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").Value = "10.47"
$ws.Range("E23").Value = "  +4.17%  "
$ws.Range("D24").Value = "4.40"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("E25").Value = "  -2.33%  "
$ws.Range("D26").Value = "161.01"
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("D27").Value = "16.91"
$ws.Range("E27").Value = "  +6.54%  "
$ws.Range("E28").Value = "  +4.77%  "
$ws.Range("D29").Value = "7.13"
$ws.Range("E29").Value = "  +5.42%  "
$ws.Range("D30").Value = "0.994"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("E31").Value = "  +10.23%  "
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("E33").Value = "  +6.20%  "
$ws.Range("E34").Value = "  +8.63%  "
$ws.Range("D35").Value = "1.572.28"
$ws.Range("E35").Value = "  +6.41%  "
$ws.Range("E36").Value = "  +5.82%  "
$ws.Range("D37").Value = "88.61"
$ws.Range("E37").Value = "  +11.18%  "
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("D39").Value = "0.629"
$ws.Range("E39").Value = "  +6.42%  "
$ws.Range("E40").Value = "  +4.60%  "
$ws.Range("D41").Value = "2.83"
$ws.Range("E41").Value = "  +6.02%  "
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("D43").Value = "0.917"
$ws.Range("E43").Value = "  +6.92%  "
$ws.Range("E44").Value = "  +5.49%  "
$ws.Range("D45").Value = "0.0520"
$ws.Range("E45").Value = "  +3.04%  "
$ws.Range("E46").Value = "  +2.85%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.947.65"
$ws.Range("E47").Value = "  +7.49%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "54.36"
$ws.Range("E48").Value = "  +4.20%  "
$ws.Range("D49").Value = "5.74"
$ws.Range("E49").Value = "  +5.27%  "
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").Value = "11.34"
$ws.Range("E51").Value = "  +21.34%  "
